# Applies the edits described by the commit:
#  - restructure sheet1 rows (blank row removal + new row added)
#  - lock the label/header cells (sheet is being protected)
#  - protect the worksheet
#  - update selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

# --- Row restructuring ------------------------------------------------
# Remove three blank rows from the gap between the "Rectangular Coupons"
# table and the "Panel" table (rows 14-17 were all blank; removing three
# of them shifts the old rows 18/19/20 up to 15/16/17).
$ws.Rows("14:16").Delete()

# Row 4 was blank before; it now holds a (still empty) formatted cell.
$ws.Range("A4").Locked = $true

# A new data row was appended under the (now) row 17 "CU-NC" row.
$ws.Range("B18").Value = 3

# --- Lock the label / header cells ------------------------------------
$ws.Range("A1").Locked = $true
$ws.Range("A2:A3").Locked = $true
$ws.Range("A5").Locked = $true
$ws.Range("A6:C6").Locked = $true
$ws.Range("A9").Locked = $true
$ws.Range("A10:E10").Locked = $true
$ws.Range("A15").Locked = $true
$ws.Range("A16:D16").Locked = $true

# --- Protect the worksheet --------------------------------------------
$ws.Protect("DE8D")

# --- Selection ----------------------------------------------------------
$ws.Range("E11").Select()
